$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.413.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.06%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.378.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.85%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.640"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.64%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.78%  "

# Row 10
$ws.Range("E10").Value = "  -4.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0919"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.35%  "

# Row 12
$ws.Range("E12").Value = "  -1.52%  "

# Row 13
$ws.Range("E13").Value = "  +0.92%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.988"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.741.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.97%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.359.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.89%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "45.400.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +21.61%  "

# Row 20
$ws.Range("E20").Value = "  -3.31%  "

# Row 21
$ws.Range("E21").Value = "  -1.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.28%  "

# Row 25
$ws.Range("E25").Value = "  -0.67%  "

# Row 26
$ws.Range("E26").Value = "  -0.08%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.49%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.41%  "

# Row 29
$ws.Range("E29").Value = "  -1.49%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0965"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.67%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.39%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "167.34"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.57%  "

# Row 34
$ws.Range("E34").Value = "  -3.66%  "

# Row 35
$ws.Range("E35").Value = "  +0.54%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.76%  "

# Row 37
$ws.Range("E37").Value = "  -4.16%  "

# Row 38
$ws.Range("E38").Value = "  +1.99%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.48%  "

# Row 40
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.58%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0356"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.19%  "

# Row 42
$ws.Range("E42").Value = "  -7.62%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.66%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.59%  "

# Row 45
$ws.Range("E45").Value = "  -4.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.75%  "

# Row 47
$ws.Range("E47").Value = "  +0.06%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.841.69"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.94%  "

# Row 49
$ws.Range("E49").Value = "  +7.22%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.67%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.46%  "
